# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 24 (pushing the previous row 24 down to row 25,
# which keeps its original values), then populate the new row 24 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 24 (and everything below it) down by one row.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly observation.
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(24, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value = 44585
$ws.Cells.Item(24, 5).Value = 13
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100102
$ws.Cells.Item(24, 8).Value = "Cítricos"
$ws.Cells.Item(24, 9).Value = 100102006
$ws.Cells.Item(24, 10).Value = "Pomelo"
$ws.Cells.Item(24, 11).Value = "Start Ruby"
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 16
$ws.Cells.Item(24, 14).Value = 200000
$ws.Cells.Item(24, 15).Value = 200000
$ws.Cells.Item(24, 16).Value = 200000
$ws.Cells.Item(24, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(24, 18).Value = "Región Metropolitana"
$ws.Cells.Item(24, 19).Value = 571
$ws.Cells.Item(24, 20).Value = 350
